$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.158.55"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -0.27%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.828.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -0.87%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'0.9980"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.47%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'242.06"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -0.97%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'0.6230"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.83%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.9988"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -0.49%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.07371"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.18%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.2914"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -1.24%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'23.17"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.87%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'0.07668"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -0.79%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.831.51"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -2.06%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'4.955"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -1.51%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'0.6675"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -1.60%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'82.57"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -0.88%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = "'0.000008976"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  -3.25%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'5.857"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -2.24%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'29.119.81"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -0.38%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'2.074.34"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -1.03%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("D20").Value = "'236.32"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.54%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'12.44"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.40%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'0.9987"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -0.56%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'7.345"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +2.18%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'0.9991"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -0.44%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'158.25"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -1.42%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'0.1417"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  +1.26%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'8.518"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.53%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'17.64"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -1.91%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = "'  -1.01%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'0.05920"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  +6.17%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = "'InternetComputer(DFINITY)"
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = "'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = "'4.080"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -2.15%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("B32").Value = "'Filecoin"
$ws.Range("B32").Style = "Normal"
$ws.Range("C32").Value = "'https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("C32").Style = "Normal"
$ws.Range("D32").Value = "'4.097"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -2.40%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'1.206"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -0.16%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.866"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +0.34%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'0.7317"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -2.75%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'1.140"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -0.96%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = "'2.601"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -2.42%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = "'2.832"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "'  +2.04%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'1.225.03"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -1.14%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("E40").Value = "'  -2.37%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'6.298"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -5.06%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.9184"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  +1.47%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'0.9986"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.44%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'102.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  -0.48%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'1.974.85"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -0.91%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'65.12"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -2.40%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'0.5047"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -1.22%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = "'  -2.40%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.4024"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.96%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'9.090"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -0.70%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("D51").Value = "'0.1130"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  +1.81%  "
$ws.Range("E51").Style = "Normal"
